# Update stats for 2025-09 (row 22 in Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C22").Value = 996
$ws.Range("D22").Value = 5846089
$ws.Range("E22").Value = 929.5737001113055
$ws.Range("G22").Value = 4.184100418410042
$ws.Range("H22").Value = 27.13421741664419
